$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
    }
}

Replace-Text "Graphic Design Institute - Employee Benefits.docx" "Graphic Design Institute での従業員手当"
Replace-Text "グラフィックデザイン研究所(GDI)は、グラフィックデザイン、Webデザイン、アニメーション、デジタルメディアに焦点を当てた商業グラフィックデザインの大手企業です。" "Graphic Design Institute (GDI) は、グラフィック デザイン、Web デザイン、アニメーション、デジタル メディアに焦点を当てた商業グラフィック デザインの大手企業です。"
Replace-Text "GDIは、成功のバックボーンである従業員に刺激的でやりがいのある職場環境を提供することに取り組んでいます。" "GDI は、成功のバックボーンである従業員に刺激的でやりがいのある職場環境を提供することに取り組んでいます。"
Replace-Text "GDIの従業員として、あなたはあなたの個人的で専門的な幸福を高める様々な利点を享受するでしょう。" "GDI の従業員は個人的かつ専門的な幸福を高める様々な利点を享受します。"
Replace-Text "GDIはあなたの健康と健康を気にし、医療、歯科、視覚の費用をカバーする包括的な医療保険プランを提供しています。" "GDI は従業員の健康とウェルネスを気にし、医療費、歯科費、眼科費をカバーする包括的な医療保険プランを提供しています。"
Replace-Text "さらに、GDI では、無料のジム メンバーシップ、フィットネス クラス、健康スクリーニング、および相談サービスを含むウェルネス プログラムを提供しています。" "さらに、GDI では、無料のジム メンバーシップ、フィットネス クラス、健康スクリーニング、相談サービスなどを含むウェルネス プログラムを提供しています。"
Replace-Text "また、GDI の従業員支援プログラム (EAP) を利用することもできます。このプログラムでは、個人および仕事に関連する問題に対する機密性の高い専門的なサポートが提供されます。" "また、GDI の従業員支援プログラム (EAP) を利用することもできます。このプログラムでは、個人および仕事関連の問題に対する機密性の高い専門的なサポートが提供されます。"
Replace-Text "退職金と節約の特典" "退職金および貯蓄給付"
Replace-Text "GDI では、将来の計画を立て、給与の最大 6% の貢献に一致する寛大な退職金プランを提供したいと考えています。" "GDI では、従業員が将来を計画できるように、給与の最大 6% までの拠出金を同額支給する手厚い退職金プランを提供しています。"
Replace-Text "さまざまな投資オプションから選択し、税繰延の成長と複合を楽しむことができます。" "さまざまな投資オプションから選択でき、課税繰り延べによる成長と複利を享受できます。"
Replace-Text "また、退職の目標と戦略を支援できる財務プランナーにもアクセスできます。" "また、退職の目標と戦略を支援する財務プランナーにもアクセスできます。"
Replace-Text "さらに、GDI には、給与控除によって給与の最大 10% を節約できる節約計画が用意されています。" "さらに、GDI には、給与控除によって給与の最大 10% を節約できる節約プランが用意されています。"
Replace-Text "GDIは、あなたの教育と開発を評価し、あなたの仕事の分野に関連する承認されたコースやプログラムのコストの最大100%をカバーする学費の払い戻しプログラムを提供しています。" "GDI では、従業員の教育と開発を重視し、従業員の仕事の分野に関連する承認されたコースやプログラムのコストの最大 100% をカバーする学費の払い戻しプログラムを提供しています。"
Replace-Text "さらに、GDI は、証明書、ボーナス、および賞を使用して実績と貢献に報酬を与える認識プログラムを提供します。" "さらに、GDI では、証明書、ボーナス、賞などを使用して、実績と貢献に報酬を与える表彰プログラムを提供します。"
Replace-Text "ワークライフバランスのメリット" "ワークライフ バランスの利点"
Replace-Text "GDI は、ワーク ライフ バランスのニーズを理解し、ニーズと好みに応じて時間と場所を調整できる柔軟な勤務スケジュールを提供します。" "GDI では、ワークライフ バランスのニーズを理解し、ニーズと好みに応じて時間や場所を調整できる柔軟な勤務スケジュールを提供します。"
Replace-Text "さらに、GDI では、休暇、病気、個人、休日の日を含む寛大な有給休暇 (PTO) ポリシーを提供します。" "さらに、GDI では、休暇、病気、個人的な日、休日を含む寛大な有給休暇 (PTO) ポリシーを提供します。"
Replace-Text "また、家族や医療上の緊急時やサバティカルなどの特別な状況に備えて、無給休暇を申請することもできます。" "また、家族の緊急時や医療上の緊急時やサバティカルなどの特別な状況に備えて、無給休暇を申請することもできます。"
Replace-Text "あなたの個性とアイデンティティを尊重し、祝う多様で包括的な文化。" "従業員の個性とアイデンティティを尊重し、祝う多様で包括的な文化。"
Replace-Text "従業員とその家族のためのイベントや活動を整理するソーシャルおよびレクリエーション プログラム。" "従業員とその家族のためのイベントや活動を企画するソーシャルおよびレクリエーション プログラム。"
Replace-Text "GDI に新しい才能をもたらすための報酬を得る紹介プログラム。" "GDI に新しい人材を紹介した従業員に報酬を与える紹介プログラム。"
Replace-Text "GDIは、あなたの生活の質とキャリアの満足度を高める包括的で競争力のある特典パッケージを提供することを誇りに思っています。" "GDI では、従業員の生活の質とキャリアの満足度を高める包括的で競争力のある特典パッケージを提供することを誇りに思っています。"
Replace-Text "GDI の従業員として、グラフィック デザインとデジタル メディアの世界にプラスの影響を与えるように努めるダイナミックでクリエイティブなチームの一員です。" "GDI の従業員は、グラフィック デザインとデジタル メディアの世界にプラスの影響を与えるように努めるダイナミックでクリエイティブなチームの一員です。"
Replace-Text "これらの利点を活用し、GDI での作業を楽しんでいただければ幸いです。" "これらの利点を活用し、GDI での勤務を楽しんでいただければ幸いです。"
Replace-Text "特典についてご質問やご不明な点がある場合は、hr@gdi.com または 555-1234 の人事部にお問い合わせください。" "特典についてご質問やご不明な点がございましたら、hr@gdi.com または 555-1234 の人事部にお問い合わせください。"

Write-Output "DONE"
